# Weekly data refresh: a new week's price observations (2 rows, Primera &
# Segunda quality) are inserted at row 972, pushing the existing rows
# 972:1079 down to 974:1081.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 972 (each Insert() shifts the
# target row, and everything below it, down by one).
$ws.Rows.Item(972).Insert()
$ws.Rows.Item(972).Insert()

# --- New row 972 : Brócoli, Primera ---
$ws.Cells.Item(972, 1).Value = 6
$ws.Cells.Item(972, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(972, 3).Value = "Metropolitana"
$ws.Cells.Item(972, 4).Value = 44946
$ws.Cells.Item(972, 5).Value = 13
$ws.Cells.Item(972, 6).Value = 100112023
$ws.Cells.Item(972, 7).Value = "Brócoli"
$ws.Cells.Item(972, 8).Value = "Sin especificar"
$ws.Cells.Item(972, 9).Value = "Primera"
$ws.Cells.Item(972, 10).Value = 12800
$ws.Cells.Item(972, 11).Value = 600
$ws.Cells.Item(972, 12).Value = 700
$ws.Cells.Item(972, 13).Value = 649
$ws.Cells.Item(972, 14).Value = "`$/unidad"
$ws.Cells.Item(972, 15).Value = "Región Metropolitana"
$ws.Cells.Item(972, 16).Value = 649
$ws.Cells.Item(972, 17).Value = 1
$ws.Cells.Item(972, 18).Value = "Hortaliza"

# --- New row 973 : Brócoli, Segunda ---
$ws.Cells.Item(973, 1).Value = 6
$ws.Cells.Item(973, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(973, 3).Value = "Metropolitana"
$ws.Cells.Item(973, 4).Value = 44946
$ws.Cells.Item(973, 5).Value = 13
$ws.Cells.Item(973, 6).Value = 100112023
$ws.Cells.Item(973, 7).Value = "Brócoli"
$ws.Cells.Item(973, 8).Value = "Sin especificar"
$ws.Cells.Item(973, 9).Value = "Segunda"
$ws.Cells.Item(973, 10).Value = 3800
$ws.Cells.Item(973, 11).Value = 400
$ws.Cells.Item(973, 12).Value = 400
$ws.Cells.Item(973, 13).Value = 400
$ws.Cells.Item(973, 14).Value = "`$/unidad"
$ws.Cells.Item(973, 15).Value = "Región Metropolitana"
$ws.Cells.Item(973, 16).Value = 400
$ws.Cells.Item(973, 17).Value = 1
$ws.Cells.Item(973, 18).Value = "Hortaliza"

# Keep the date cells formatted the same way as the rest of column D.
$ws.Range("D972:D973").NumberFormat = $ws.Range("D974").NumberFormat
